# Add the "新規媒体コード" (new media code) values for the remaining rows of
# column F (rows 46-94) that did not already carry an XLOOKUP-derived value.
# Source workbook: 媒体コードマスタ.xlsx
#
# Note on scope: this script reproduces every change that is reachable
# through the Excel object model exposed here (new literal string values in
# F46:F94, matching shared-string growth, matching cell style/row spans).
# Two items from the upstream diff are pure Excel-internal bookkeeping that
# have no COM-exposed setter and are not touched by real user actions
# either: the x15ac:absPath cache (the folder Excel was opened from) and
# the xr:revisionPtr documentId GUID (a collaboration/session id). Both are
# regenerated by Excel itself on save, not something a script sets. The
# orphaned (unreferenced) dxf in styles.xml is likewise left as-is: the
# workbook model only supports *adding*/*updating* dxfs, there is no
# delete/prune primitive anywhere in the object model to remove an unused
# one, so it cannot be cleared without risking corrupting/duplicating it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing formula-bearing F-column cell (border + font)
# onto the new F46:F94 range before writing values, so the added cells keep
# the same cell style (s="3") as the rest of column F instead of defaulting
# to an unstyled cell.
$ws.Range("F44").Copy()
$ws.Range("F46:F94").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New 媒体コード (media code) values for rows 46-94, column F —
# these rows previously had no XLOOKUP match against the external
# workbook, so the values are entered directly as literal text.
$ws.Range("F46").Value = "swpbh"
$ws.Range("F47").Value = "swpjs"
$ws.Range("F48").Value = "swpau"
$ws.Range("F49").Value = "swpzc"
$ws.Range("F50").Value = "swp3f"
$ws.Range("F51").Value = "swp6j"
$ws.Range("F52").Value = "swpev"
$ws.Range("F53").Value = "swpm9"
$ws.Range("F54").Value = "swpq1"
$ws.Range("F55").Value = "swpvx"
$ws.Range("F56").Value = "swp2m"
$ws.Range("F57").Value = "swpij"
$ws.Range("F58").Value = "swp5d"
$ws.Range("F59").Value = "swpzz"
$ws.Range("F60").Value = "swpe1"
$ws.Range("F61").Value = "swpn9"
$ws.Range("F62").Value = "swpvj"
$ws.Range("F63").Value = "swpwz"
$ws.Range("F64").Value = "swpiz"
$ws.Range("F65").Value = "swpgi"
$ws.Range("F66").Value = "swpvp"
$ws.Range("F67").Value = "swpjo"
$ws.Range("F68").Value = "swpcg"
$ws.Range("F69").Value = "swpo7"
$ws.Range("F70").Value = "swpeh"
$ws.Range("F71").Value = "swpnp"
$ws.Range("F72").Value = "swppz"
$ws.Range("F73").Value = "swpm6"
$ws.Range("F74").Value = "swpor"
$ws.Range("F75").Value = "swpku"
$ws.Range("F76").Value = "swpo2"
$ws.Range("F77").Value = "swpx3"
$ws.Range("F78").Value = "swpdy"
$ws.Range("F79").Value = "swpjj"
$ws.Range("F80").Value = "swpyv"
$ws.Range("F81").Value = "swpba"
$ws.Range("F82").Value = "swpco"
$ws.Range("F83").Value = "swpmv"
$ws.Range("F84").Value = "swpaa"
$ws.Range("F85").Value = "swp6m"
$ws.Range("F86").Value = "swp0a"
$ws.Range("F87").Value = "swpvo"
$ws.Range("F88").Value = "swpj6"
$ws.Range("F89").Value = "swpu5"
$ws.Range("F90").Value = "swpzx"
$ws.Range("F91").Value = "swpx0"
$ws.Range("F92").Value = "swpz8"
$ws.Range("F93").Value = "swphr"
$ws.Range("F94").Value = "swp9o"
